$wb = $excel.ActiveWorkbook

$wsPIR = $wb.Worksheets.Item("PIR")
$PIRRows = @(
    @("308", "2026-01-28", "12:30:08", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("309", "2026-01-28", "12:30:09", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("310", "2026-01-28", "12:30:14", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("311", "2026-01-28", "12:30:19", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("312", "2026-01-28", "12:30:26", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("313", "2026-01-28", "12:30:29", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("314", "2026-01-28", "12:30:34", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("315", "2026-01-28", "12:30:39", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("316", "2026-01-28", "12:30:46", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("317", "2026-01-28", "12:30:50", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("318", "2026-01-28", "12:30:54", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("319", "2026-01-28", "12:30:59", "12:00", "Bathroom", "No Motion", "Inactive"),
    @("320", "2026-01-28", "12:31:04", "12:00", "Bathroom", "No Motion", "Inactive")
)
foreach ($r in $PIRRows) {
    $rowNum = $r[0]
    for ($col = 1; $col -le 6; $col++) {
        $cell = $wsPIR.Cells.Item($rowNum, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $r[$col]
        $cell.Style = "Normal"
    }
}

$wsHumidity = $wb.Worksheets.Item("Humidity")
$HumidityRows = @(
    @("286", "2026-01-28", "12:30:05", "12:00", "Bathroom", "87.6%", "Active"),
    @("287", "2026-01-28", "12:30:06", "12:00", "Bathroom", "86.7%", "Active"),
    @("288", "2026-01-28", "12:30:16", "12:00", "Bathroom", "86.8%", "Active"),
    @("289", "2026-01-28", "12:30:20", "12:00", "Bathroom", "87.7%", "Active"),
    @("290", "2026-01-28", "12:30:24", "12:00", "Bathroom", "86.7%", "Active"),
    @("291", "2026-01-28", "12:30:28", "12:00", "Bathroom", "87.6%", "Active"),
    @("292", "2026-01-28", "12:30:32", "12:00", "Bathroom", "87.6%", "Active"),
    @("293", "2026-01-28", "12:30:36", "12:00", "Bathroom", "86.7%", "Active"),
    @("294", "2026-01-28", "12:30:40", "12:00", "Bathroom", "87.6%", "Active"),
    @("295", "2026-01-28", "12:30:44", "12:00", "Bathroom", "86.7%", "Active"),
    @("296", "2026-01-28", "12:30:48", "12:00", "Bathroom", "87.6%", "Active"),
    @("297", "2026-01-28", "12:30:52", "12:00", "Bathroom", "87.6%", "Active"),
    @("298", "2026-01-28", "12:31:00", "12:00", "Bathroom", "87.6%", "Active")
)
foreach ($r in $HumidityRows) {
    $rowNum = $r[0]
    for ($col = 1; $col -le 6; $col++) {
        $cell = $wsHumidity.Cells.Item($rowNum, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $r[$col]
        $cell.Style = "Normal"
    }
}

$wsTemperature = $wb.Worksheets.Item("Temperature")
$TemperatureRows = @(
    @("286", "2026-01-28", "12:30:05", "12:00", "Bathroom", "22.9C", "Active"),
    @("287", "2026-01-28", "12:30:07", "12:00", "Bathroom", "22.9C", "Active"),
    @("288", "2026-01-28", "12:30:17", "12:00", "Bathroom", "22.9C", "Active"),
    @("289", "2026-01-28", "12:30:21", "12:00", "Bathroom", "22.9C", "Active"),
    @("290", "2026-01-28", "12:30:25", "12:00", "Bathroom", "22.9C", "Active"),
    @("291", "2026-01-28", "12:30:29", "12:00", "Bathroom", "22.9C", "Active"),
    @("292", "2026-01-28", "12:30:33", "12:00", "Bathroom", "22.9C", "Active"),
    @("293", "2026-01-28", "12:30:37", "12:00", "Bathroom", "22.9C", "Active"),
    @("294", "2026-01-28", "12:30:41", "12:00", "Bathroom", "22.9C", "Active"),
    @("295", "2026-01-28", "12:30:45", "12:00", "Bathroom", "22.9C", "Active"),
    @("296", "2026-01-28", "12:30:49", "12:00", "Bathroom", "22.9C", "Active"),
    @("297", "2026-01-28", "12:30:53", "12:00", "Bathroom", "22.9C", "Active"),
    @("298", "2026-01-28", "12:31:01", "12:00", "Bathroom", "22.9C", "Active")
)
foreach ($r in $TemperatureRows) {
    $rowNum = $r[0]
    for ($col = 1; $col -le 6; $col++) {
        $cell = $wsTemperature.Cells.Item($rowNum, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $r[$col]
        $cell.Style = "Normal"
    }
}
